$wb = $excel.ActiveWorkbook

# --- "LoginLogoutRegistration" sheet (physical sheet2.xml) gets a new
#     "version" column inserted before the existing "firstname" column. ---
$ws2 = $wb.Worksheets.Item("LoginLogoutRegistration")

$ws2.Columns("E:E").Insert()

$ws2.Range("E1").Value = "'version"
$ws2.Range("E2:E7").Value = "'120.0"

$ws2.Range("K1").Select()

# --- "Run Manager" sheet (physical sheet1.xml) just had its selection
#     changed to B3. ---
$ws1 = $wb.Worksheets.Item("Run Manager")
$ws1.Range("B3").Select()

$ws2.Activate()
